$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Bomba Dosadora com Diafragma" nominal power value (D10)
$ws.Range("D10").Value = 0.45

# Recalculate so dependent formula cells (F10, F11) pick up the new value
$excel.Calculate()

# Reflect the new selection used while reviewing the updated rows
$ws.Range("F6:F10").Select()
